$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells so numeric-looking strings
# (e.g. "57.638.52", "1.00") are preserved exactly as text, not converted
# to numbers/dates by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values from the diff
$ws.Range("D2").Value = '57.638.52'
$ws.Range("E2").Value = '  +0.28%  '
$ws.Range("D3").Value = '3.118.87'
$ws.Range("E3").Value = '  +0.22%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '531.40'
$ws.Range("E5").Value = '  +1.14%  '
$ws.Range("D6").Value = '138.12'
$ws.Range("E6").Value = '  +0.83%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.114.93'
$ws.Range("E8").Value = '  +0.12%  '
$ws.Range("D9").Value = '0.469'
$ws.Range("E9").Value = '  +4.96%  '
$ws.Range("E10").Value = '  +0.68%  '
$ws.Range("E11").Value = '  +0.32%  '
$ws.Range("E12").Value = '  +4.47%  '
$ws.Range("E13").Value = '  +1.47%  '
$ws.Range("D14").Value = '3.652.81'
$ws.Range("E14").Value = '  +0.04%  '
$ws.Range("D15").Value = '25.67'
$ws.Range("E15").Value = '  +1.70%  '
$ws.Range("D16").Value = '0.0000164'
$ws.Range("E16").Value = '  +1.03%  '
$ws.Range("D17").Value = '57.767.71'
$ws.Range("E17").Value = '  +0.33%  '
$ws.Range("D18").Value = '3.113.90'
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("E19").Value = '  +1.60%  '
$ws.Range("D20").Value = '12.68'
$ws.Range("E20").Value = '  +2.37%  '
$ws.Range("E21").Value = '  +2.53%  '
$ws.Range("D22").Value = '361.99'
$ws.Range("E22").Value = '  +4.21%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").Value = '69.04'
$ws.Range("E24").Value = '  +1.64%  '
$ws.Range("D25").Value = '0.505'
$ws.Range("E25").Value = '  +0.83%  '
$ws.Range("D26").Value = '0.167'
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("E27").Value = '  +0.09%  '
$ws.Range("D28").Value = '0.0₃0863'
$ws.Range("E28").Value = '  -4.43%  '
$ws.Range("D29").Value = '7.31'
$ws.Range("E29").Value = '  -1.38%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '1.87'
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("B31").Value = 'RenderToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D31").Value = '6.08'
$ws.Range("E31").Value = '  +0.60%  '
$ws.Range("D32").Value = '21.38'
$ws.Range("E32").Value = '  +2.05%  '
$ws.Range("D33").Value = '5.11'
$ws.Range("E33").Value = '  +3.23%  '
$ws.Range("D34").Value = '159.62'
$ws.Range("E34").Value = '  +0.63%  '
$ws.Range("E35").Value = '  -1.37%  '
$ws.Range("D36").Value = '6.05'
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '1.28'
$ws.Range("E37").Value = '  +3.66%  '
$ws.Range("B38").Value = 'EnergySwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D38").Value = '25.52'
$ws.Range("E38").Value = '  -1.76%  '
$ws.Range("E39").Value = '  +2.63%  '
$ws.Range("D40").Value = '0.0670'
$ws.Range("E40").Value = '  +1.09%  '
$ws.Range("D41").Value = '2.497.43'
$ws.Range("E41").Value = '  +5.85%  '
$ws.Range("D42").Value = '4.01'
$ws.Range("E42").Value = '  -4.39%  '
$ws.Range("D43").Value = '0.696'
$ws.Range("E43").Value = '  -0.71%  '
$ws.Range("D44").Value = '37.73'
$ws.Range("E44").Value = '  +3.39%  '
$ws.Range("D45").Value = '0.0270'
$ws.Range("E45").Value = '  +1.36%  '
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  -0.05%  '
$ws.Range("B47").Value = 'ONDO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D47").Value = '0.980'
$ws.Range("E47").Value = '  +2.15%  '
$ws.Range("B48").Value = 'Cosmos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D48").Value = '6.07'
$ws.Range("E48").Value = '  +0.91%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '19.66'
$ws.Range("E49").Value = '  -1.29%  '
$ws.Range("B50").Value = 'SuiNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D50").Value = '0.736'
$ws.Range("E50").Value = '  -3.47%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").Value = '0.0909'
$ws.Range("E51").Value = '  +1.91%  '

# Restore default (Normal) style on the Price cells so no stray
# number-format style is left attached to them.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
